$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last_edited_time" (column D) values to reflect the new
# timestamps recorded when bonus/penalty totals were consolidated
# into the personal report.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-08-03T21:27:00.000Z"
}

for ($r = 8; $r -le 22; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-08-03T21:28:00.000Z"
}
